$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 120
$ws1.Range("F4").Value = 437
$ws1.Range("F7").Value = 1211
$ws1.Range("F8").Value = 406
$ws1.Range("F12").Value = 383
$ws1.Range("F14").Value = 799
$ws1.Range("F17").Value = 294
$ws1.Range("F19").Value = 1026
$ws1.Range("F22").Value = 90
$ws1.Range("F23").Value = 389

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 348

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 348
$ws4.Range("F4").Value = 120
$ws4.Range("F6").Value = 437
$ws4.Range("F9").Value = 1211
$ws4.Range("F10").Value = 406
$ws4.Range("F17").Value = 383
$ws4.Range("F21").Value = 799
$ws4.Range("F24").Value = 294
$ws4.Range("F26").Value = 1026
$ws4.Range("F31").Value = 90
$ws4.Range("F32").Value = 389
